# This sheet holds a weekly price series for "Poroto verde" at the
# "Feria Lagunitas de Puerto Montt" market, most-recent-first. A new
# weekly observation (2022-08-26, one week after the former top row's
# 2022-08-19) is being added at the top of the data block (row 44, just
# below the header row and the very first historical record in row 2-43
# stays put). All the existing observations that used to live in rows
# 44-95 shift down by one row, to 45-96, and the sheet's used range grows
# from R95 to R96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 44; this pushes the current rows 44-95 down
# to 45-96 (Excel's normal "insert row above" behaviour), growing the
# sheet from A1:R95 to A1:R96.
$ws.Rows.Item(44).Insert()

# Fill the newly inserted row 44 with the new weekly record. It repeats
# the same market/product/quality/price-unit/origin/unit-price metadata
# as the (now shifted-down) row that used to be on top, only the date
# and the traded volume differ.
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44799
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112031
$ws.Range("G44").Value = "Poroto verde"
$ws.Range("H44").Value = "Magnum"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = 35000
$ws.Range("N44").Value = "$/malla 25 kilos"
$ws.Range("O44").Value = "Perú"
$ws.Range("P44").Value = 1400
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
